$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

$xlPasteFormats = -4122

function CopyFormat($srcAddr, $dstAddr) {
  $ws.Range($srcAddr).Copy() | Out-Null
  $ws.Range($dstAddr).PasteSpecial($xlPasteFormats) | Out-Null
}

# --- 1) Grow the table by 9 rows. ListRows.Add(112) repeatedly inserts a
#     fresh blank row immediately above the current Totals row (pushing the
#     Totals row down by one each time), so after 9 calls the table spans
#     A1:E122 with the Totals row sitting at 122 and 9 new blank data rows
#     at 113-121. ---
for ($i = 0; $i -lt 9; $i++) {
  $lo.ListRows.Add(112) | Out-Null
}
$excel.CutCopyMode = 0

# --- 2) Relocate the grand-total row's formatting + content from its old
#     home (row 113) down to its new home (row 122), column by column so we
#     never touch B/D (which the totals row doesn't use). ---
CopyFormat "A113" "A122"
CopyFormat "C113" "C122"
CopyFormat "E113" "E122"

$totalLabel = $ws.Range("A113").Value()
$ws.Range("A113").ClearContents()
$ws.Range("C113").ClearContents()
$ws.Range("E113").ClearContents()

$ws.Range("A122").Value = $totalLabel
$ws.Range("C122").Formula = "=SUBTOTAL(109,Table1[Hours])"

# --- 3) Row 113 becomes the new "Week 14" section header, styled like the
#     previous week-header rows (e.g. row 104). ---
CopyFormat "A104" "A113"
CopyFormat "C104" "C113"
CopyFormat "E104" "E113"
$ws.Range("A113").Value = "Week 14"

# --- 4) Row 114: first (and only, so far) logged day of the new week. ---
CopyFormat "A105" "A114"
CopyFormat "B105" "B114"
CopyFormat "C105" "C114"
CopyFormat "D105" "D114"
CopyFormat "E105" "E114"

$ws.Range("A114").Value = 44535
$ws.Range("B114").Value = "JS101: Programming Foundations with JavaScript"
$ws.Range("C114").Value = 1
$ws.Range("D114").Value = "6.9 convert data to one cards object, displayHands function"

# --- 5) Rows 115-120: remaining days of the week - date only, matching the
#     plain date-row style (e.g. row 106). ---
$dates = @(44536, 44537, 44538, 44539, 44540, 44541)
$r = 115
foreach ($d in $dates) {
  CopyFormat "A106" "A$r"
  CopyFormat "C106" "C$r"
  CopyFormat "E106" "E$r"
  $ws.Range("A$r").Value = $d
  $r = $r + 1
}

# --- 6) Row 121: "Weekly Total" row for the new week, styled like the
#     previous week's weekly-total row (112). ---
CopyFormat "A112" "A121"
CopyFormat "B112" "B121"
CopyFormat "C112" "C121"
CopyFormat "D112" "D121"
CopyFormat "E112" "E121"

$ws.Range("B121").Value = "Weekly Total"
$ws.Range("D121").Formula = "=SUM(C114:C120)"

$excel.CutCopyMode = 0

# --- 7) Sheet view / selection tweaks from the diff. ---
$activeWindow = $excel.ActiveWindow
$activeWindow.ScrollRow = 99
$activeWindow.ScrollColumn = 3
$ws.Range("D114").Select() | Out-Null
